$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: merge a run with the immediately preceding run (same formatting)
# by re-"typing" its current full text via Find/Replace scoped to a Range.
# Word's replace engine coalesces the new text into the preceding run when
# the formatting matches, exactly mirroring what happens when a user
# retypes/overwrites text that spans a run boundary.
# ---------------------------------------------------------------------------

function Merge-Text {
    param($range, [string]$text)
    $null = $range.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 2)
}

# ---------------------------------------------------------------------------
# 1) "- codifica degli script javascript della galleria" paragraph:
#    merge " " + "codifica degli script " -> " codifica degli script "
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(26)
Merge-Text $p.Range "codifica degli script "

# ---------------------------------------------------------------------------
# 2) "- codifica e manutenzione degli script PERL/CGI pagine amministratore"
#    merge "pagine " + "amministratore" -> " pagine amministratore"
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(29)
Merge-Text $p.Range "amministratore"

# ---------------------------------------------------------------------------
# 3) "- codifica e manutenzione di codice CSS e XHTML"
#    merge "codice C" + "SS e XHTML" -> "codice CSS e XHTML"
#    then merge " " + "codifica e manutenzione di codice CSS e XHTML"
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(30)
Merge-Text $p.Range "SS e XHTML"
$p = $d.Paragraphs(30)
Merge-Text $p.Range "codifica e manutenzione di codice CSS e XHTML"

# ---------------------------------------------------------------------------
# 4) "- codifica e manutenzione degli script PERL/CGI pagina tariffe"
#    merge " " + "tariffe" -> " tariffe"
#    then merge " pagina" + " tariffe" -> " pagina tariffe"
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(36)
Merge-Text $p.Range "tariffe"
$p = $d.Paragraphs(36)
Merge-Text $p.Range " tariffe"

# ---------------------------------------------------------------------------
# 5) "- codifica e manutenzione di codice CSS e XHTML" (second occurrence)
#    merge "codifica e manutenzione " + "di codice CSS e XHTML"
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(37)
Merge-Text $p.Range "di codice CSS e XHTML"

# ---------------------------------------------------------------------------
# 6) "//DOMANDA CRISTIAN: la parte admin è ottimizzata per il mobile?"
#    color the whole paragraph red, then split "mobile?" into "mobi" | "le?"
#    with the _GoBack bookmark positioned between them (Bookmarks.Add moves
#    any existing bookmark of the same name, so the old one is removed
#    automatically).
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(67)
$p.Range.Font.Color = 255

$rng = $d.Content
$null = $rng.Find.Execute("mobile?")
$bmRange = $d.Range($rng.Start + 4, $rng.Start + 4)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 7) "Pagina vista da chi soffre di Deuteranopia"
#    merge "Pagina vista da chi soffre di" + " " + "Deuteranopia"
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(77)
Merge-Text $p.Range "Deuteranopia"
$p = $d.Paragraphs(77)
Merge-Text $p.Range " Deuteranopia"

# ---------------------------------------------------------------------------
# 8) Browser compatibility list merges
# ---------------------------------------------------------------------------
# "- Microsoft Edge versione 39.14915.1000.0"
$p = $d.Paragraphs(93)
Merge-Text $p.Range " 39.14915.1000.0"
$p = $d.Paragraphs(93)
Merge-Text $p.Range "versione 39.14915.1000.0"

# "- Opera versione 39.0.2256.48"
$p = $d.Paragraphs(94)
Merge-Text $p.Range " 39.0.2256.48"
$p = $d.Paragraphs(94)
Merge-Text $p.Range "versione 39.0.2256.48"

# "- Internet Explorer versione 7,8,9"
$p = $d.Paragraphs(95)
Merge-Text $p.Range " 7,8,9"
$p = $d.Paragraphs(95)
Merge-Text $p.Range "versione 7,8,9"

# "- Safari versione mobile"
$p = $d.Paragraphs(96)
Merge-Text $p.Range " mobile"
$p = $d.Paragraphs(96)
Merge-Text $p.Range "versione mobile"

# "- Mozilla Firefox versione //VERSIONE FIREFOX"
$p = $d.Paragraphs(97)
Merge-Text $p.Range "versione "
$p = $d.Paragraphs(97)
Merge-Text $p.Range " versione "

Write-Output "done"
